$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "ID hráče" in F2
$ws.Range("F2").Value = "ID hráče"

# Match F2's formatting (border) to the rest of the header row before re-aligning
$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Center the whole header row (A2:F2), which also centers the pre-existing cells
$ws.Range("A2:F2").HorizontalAlignment = -4108  # xlCenter

# Widen column F to match the other data columns
$ws.Columns.Item(6).ColumnWidth = 17

# Update selection to reflect the new active cell
$ws.Range("G2").Select()
